$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 36 (pushes the old "keT3 (30)" block, and
# everything below it, down by two rows). Ordinary formulas that use
# plain cell references (e.g. "=F44") are shifted automatically by the
# Insert, but workbook-level defined names are NOT - those are fixed up
# by hand below, same as the real ABCdrug/ABCsol/eps/k13D/k13DS/k13S/
# keT3_/P/Rcap/Rkrogh shifts seen in the target diff.
$ws.Rows.Item(36).Resize(2).Insert()

$shiftedNames = @{
    "ABCdrug" = 40
    "ABCsol"  = 41
    "eps"     = 45
    "k13D"    = 46
    "k13DS"   = 48
    "k13S"    = 47
    "keT3_"   = 38
    "P"       = 42
    "Rcap"    = 43
    "Rkrogh"  = 44
}
foreach ($name in $shiftedNames.Keys) {
    $row = $shiftedNames[$name]
    $wb.Names.Item($name).RefersTo = "=Sheet1!`$F`$$row"
}

# New row 36: keT1, computed from the (now-shifted) keT3_ named range.
$ws.Range("E36").Value = "keT1"
$ws.Range("F36").Formula = "=keT3_"
$ws.Range("G36").Value = "1/d"
$ws.Range("H36").Value = "calc"

# New row 37: keDT1, simply mirrors the new F36 value.
$ws.Range("E37").Value = "keDT1"
$ws.Range("F37").Formula = "=F36"
$ws.Range("G37").Value = "1/d"
$ws.Range("H37").Value = "calc"

# Sheet view: scrolled down a bit, with D25 selected.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D25").Select()
